$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.931.92"
$ws.Range("E2").Value = "'  +0.44%  "
$ws.Range("D3").Value = "'2.915.66"
$ws.Range("E3").Value = "'  -0.35%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'588.04"
$ws.Range("E5").Value = "'  -1.26%  "
$ws.Range("D6").Value = "'146.04"
$ws.Range("E6").Value = "'  +3.16%  "
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E8").Value = "'  +1.42%  "
$ws.Range("D9").Value = "'2.914.16"
$ws.Range("E9").Value = "'  -0.34%  "
$ws.Range("D10").Value = "'7.00"
$ws.Range("E10").Value = "'  -2.65%  "
$ws.Range("E11").Value = "'  +6.86%  "
$ws.Range("E12").Value = "'  -0.99%  "
$ws.Range("E13").Value = "'  +6.72%  "
$ws.Range("E14").Value = "'  -1.69%  "
$ws.Range("E15").Value = "'  -1.30%  "
$ws.Range("D16").Value = "'3.400.47"
$ws.Range("E16").Value = "'  -0.32%  "
$ws.Range("D17").Value = "'61.940.78"
$ws.Range("E17").Value = "'  +0.63%  "
$ws.Range("D18").Value = "'6.60"
$ws.Range("E18").Value = "'  -0.56%  "
$ws.Range("D19").Value = "'2.918.72"
$ws.Range("E19").Value = "'  +0.06%  "
$ws.Range("D20").Value = "'434.70"
$ws.Range("E20").Value = "'  +0.28%  "
$ws.Range("D21").Value = "'13.42"
$ws.Range("E21").Value = "'  -0.31%  "
$ws.Range("D22").Value = "'0.660"
$ws.Range("E22").Value = "'  -1.41%  "
$ws.Range("E23").Value = "'  -1.41%  "
$ws.Range("D24").Value = "'80.93"
$ws.Range("E24").Value = "'  -0.11%  "
$ws.Range("D25").Value = "'10.95"
$ws.Range("E25").Value = "'  +2.82%  "
$ws.Range("D26").Value = "'11.90"
$ws.Range("E26").Value = "'  +1.95%  "
$ws.Range("D27").Value = "'2.09"
$ws.Range("E27").Value = "'  -1.00%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "'  -0.09%  "
$ws.Range("D29").Value = "'7.32"
$ws.Range("E29").Value = "'  +6.97%  "
$ws.Range("E30").Value = "'  +21.82%  "
$ws.Range("D31").Value = "'2.57"
$ws.Range("E31").Value = "'  -0.59%  "
$ws.Range("D32").Value = "'2.11"
$ws.Range("E32").Value = "'  +0.99%  "
$ws.Range("E33").Value = "'  +3.75%  "
$ws.Range("D34").Value = "'26.05"
$ws.Range("E34").Value = "'  -0.36%  "
$ws.Range("E35").Value = "'  -0.09%  "
$ws.Range("D36").Value = "'0.977"
$ws.Range("E36").Value = "'  -0.86%  "
$ws.Range("E37").Value = "'  +8.72%  "
$ws.Range("D38").Value = "'5.53"
$ws.Range("E38").Value = "'  -0.51%  "
$ws.Range("D39").Value = "'49.29"
$ws.Range("E39").Value = "'  +0.18%  "
$ws.Range("E40").Value = "'  +2.99%  "
$ws.Range("D41").Value = "'8.37"
$ws.Range("E41").Value = "'  -1.33%  "
$ws.Range("E42").Value = "'  -1.52%  "
$ws.Range("D43").Value = "'0.274"
$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("D44").Value = "'38.93"
$ws.Range("E44").Value = "'  +1.10%  "
$ws.Range("D45").Value = "'2.699.64"
$ws.Range("E45").Value = "'  +0.61%  "
$ws.Range("D46").Value = "'134.97"
$ws.Range("E46").Value = "'  +1.12%  "
$ws.Range("D47").Value = "'0.0339"
$ws.Range("E47").Value = "'  +0.76%  "
$ws.Range("D48").Value = "'346.66"
$ws.Range("E48").Value = "'  -2.81%  "
$ws.Range("E50").Value = "'  +0.62%  "
$ws.Range("D51").Value = "'22.53"
$ws.Range("E51").Value = "'  -0.76%  "
